$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "read" method to "inscription" (cell A2).
$ws.Range("A2").Value = "inscription"

# Move the active selection to A2.
$ws.Range("A2").Select()
